# Auto-generated edit script: refresh market-data derived profit columns
# across multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 2077.5
$ws.Range("I15").Value = 2077.5
$ws.Range("K15").Value = 6232.5
$ws.Range("M15").Value = -6063.5

# Row 28
$ws.Range("H28").Value = 11793.182
$ws.Range("I28").Value = 2056.75
$ws.Range("K28").Value = 2056.75
$ws.Range("M28").Value = -1571.75

# Row 92
$ws.Range("H92").Value = 1132.125
$ws.Range("I92").Value = 117.833336
$ws.Range("J92").Value = 4175
$ws.Range("K92").Value = 117.833336
$ws.Range("L92").Value = 4175
$ws.Range("M92").Value = 1130.166664
$ws.Range("N92").Value = -6671

# Row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# Row 138
$ws.Range("H138").Value = 2705.0625
$ws.Range("J138").Value = 3339.2222
$ws.Range("L138").Value = 10017.6666
$ws.Range("N138").Value = -20297.6666

$ws = $wb.Worksheets.Item("ARM")
# Row 110
$ws.Range("H110").Value = 3216.6428
$ws.Range("I110").Value = 2942
$ws.Range("J110").Value = 3903.25
$ws.Range("K110").Value = 2942
$ws.Range("L110").Value = 3903.25
$ws.Range("M110").Value = -897
$ws.Range("N110").Value = -7993.25

# Row 122
$ws.Range("H122").Value = 4972.6665
$ws.Range("I122").Value = 4966
$ws.Range("J122").Value = 4996
$ws.Range("K122").Value = 14898
$ws.Range("L122").Value = 14988
$ws.Range("M122").Value = -12448
$ws.Range("N122").Value = -19888

# Row 132
$ws.Range("H132").Value = 3414.0908
$ws.Range("I132").Value = 1882.375
$ws.Range("J132").Value = 7498.6665
$ws.Range("K132").Value = 5647.125
$ws.Range("L132").Value = 22495.9995
$ws.Range("M132").Value = -3117.125
$ws.Range("N132").Value = -27555.9995

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3866.4119
$ws.Range("I86").Value = 1626.3334
$ws.Range("K86").Value = 1626.3334
$ws.Range("M86").Value = -503.3334

# Row 89
$ws.Range("H89").Value = 3866.4119
$ws.Range("I89").Value = 1626.3334
$ws.Range("K89").Value = 8131.666999999999
$ws.Range("M89").Value = -2515.666999999999

# Row 94
$ws.Range("H94").Value = 893.75
$ws.Range("I94").Value = 860
$ws.Range("K94").Value = 860
$ws.Range("M94").Value = -409

# Row 107
$ws.Range("H107").Value = 3428.9666
$ws.Range("I107").Value = 1093.5
$ws.Range("K107").Value = 1093.5
$ws.Range("M107").Value = 826.5

# Row 134
$ws.Range("H134").Value = 3543.7778
$ws.Range("J134").Value = 4665
$ws.Range("L134").Value = 13995
$ws.Range("N134").Value = -19065

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1413.4762
$ws.Range("I22").Value = 865.3333
$ws.Range("J22").Value = 2144.3333
$ws.Range("K22").Value = 865.3333
$ws.Range("L22").Value = 2144.3333
$ws.Range("M22").Value = -515.3333
$ws.Range("N22").Value = -2844.3333

# Row 94
$ws.Range("H94").Value = 3292.375
$ws.Range("I94").Value = 2666.2
$ws.Range("J94").Value = 3577
$ws.Range("K94").Value = 2666.2
$ws.Range("L94").Value = 3577
$ws.Range("M94").Value = -2215.2
$ws.Range("N94").Value = -4479

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 2619.5
$ws.Range("I5").Value = 2059.6667
$ws.Range("J5").Value = 4299
$ws.Range("K5").Value = 6179.000100000001
$ws.Range("L5").Value = 12897
$ws.Range("M5").Value = -6067.000100000001
$ws.Range("N5").Value = -13121

# Row 135
$ws.Range("H135").Value = 2619.5
$ws.Range("I135").Value = 2059.6667
$ws.Range("J135").Value = 4299
$ws.Range("K135").Value = 18537.0003
$ws.Range("L135").Value = 38691
$ws.Range("M135").Value = -16002.0003
$ws.Range("N135").Value = -43761

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 8998.143
$ws.Range("I70").Value = 3745.5
$ws.Range("J70").Value = 16001.667
$ws.Range("K70").Value = 3745.5
$ws.Range("L70").Value = 16001.667
$ws.Range("M70").Value = -3475.5
$ws.Range("N70").Value = -16541.667

# Row 73
$ws.Range("H73").Value = 8998.143
$ws.Range("I73").Value = 3745.5
$ws.Range("J73").Value = 16001.667
$ws.Range("K73").Value = 3745.5
$ws.Range("L73").Value = 16001.667
$ws.Range("M73").Value = -2809.5
$ws.Range("N73").Value = -17873.667

# Row 80
$ws.Range("H80").Value = 1829.8
$ws.Range("I80").Value = 1813.1428
$ws.Range("J80").Value = 1868.6666
$ws.Range("K80").Value = 1813.1428
$ws.Range("L80").Value = 1868.6666
$ws.Range("M80").Value = -815.1428000000001
$ws.Range("N80").Value = -3864.6666

# Row 83
$ws.Range("H83").Value = 1829.8
$ws.Range("I83").Value = 1813.1428
$ws.Range("J83").Value = 1868.6666
$ws.Range("K83").Value = 9065.714
$ws.Range("L83").Value = 9343.333000000001
$ws.Range("M83").Value = -4073.714
$ws.Range("N83").Value = -19327.333

# Row 123
$ws.Range("H123").Value = 109999
$ws.Range("J123").Value = 109999
$ws.Range("L123").Value = 109999
$ws.Range("N123").Value = -114899

# Row 132
$ws.Range("H132").Value = 29513.824
$ws.Range("I132").Value = 37323.535
$ws.Range("K132").Value = 111970.605
$ws.Range("M132").Value = -109440.605

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1017.3333
$ws.Range("I22").Value = 1050.5
$ws.Range("J22").Value = 951
$ws.Range("K22").Value = 1050.5
$ws.Range("L22").Value = 951
$ws.Range("M22").Value = -755.5
$ws.Range("N22").Value = -1541

# Row 27
$ws.Range("H27").Value = 1017.3333
$ws.Range("I27").Value = 1050.5
$ws.Range("J27").Value = 951
$ws.Range("K27").Value = 1050.5
$ws.Range("L27").Value = 951
$ws.Range("M27").Value = -943.5
$ws.Range("N27").Value = -1165

# Row 46
$ws.Range("H46").Value = 4574.2856
$ws.Range("I46").Value = 3891.111
$ws.Range("J46").Value = 4897.8945
$ws.Range("K46").Value = 3891.111
$ws.Range("L46").Value = 4897.8945
$ws.Range("M46").Value = -3703.111
$ws.Range("N46").Value = -5273.8945

# Row 122
$ws.Range("H122").Value = 5152.6924
$ws.Range("I122").Value = 4999.7144
$ws.Range("J122").Value = 5331.1665
$ws.Range("K122").Value = 14999.1432
$ws.Range("L122").Value = 15993.4995
$ws.Range("M122").Value = -12549.1432
$ws.Range("N122").Value = -20893.4995

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 3382.9412
$ws.Range("I132").Value = 1967.125
$ws.Range("K132").Value = 5901.375
$ws.Range("M132").Value = -3371.375

# Row 138
$ws.Range("H138").Value = 95000
$ws.Range("J138").Value = 95000
$ws.Range("L138").Value = 95000
$ws.Range("N138").Value = -105280
